$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New pretty-printed JSON content (replaces old compact Python-repr string).
# Single-quoted here-string: no interpolation/escaping, fully literal.
$newText = @'
questions = [
    {
        "title": "You want to initialize a repository in an existing directory. Which command achieves this?",
        "ques_type": 2,
        "options": [
            "git clone",
            "git init",
            "git create",
            "git fetch"
        ],
        "score": "git init"
    },
    {
        "title": "You want to create and switch to a new branch called \u201cnew-branches\u201d.  Which command achieves this?",
        "ques_type": 2,
        "options": [
            "git branch -D new-branches ",
            "git checkout new-branches",
            "git checkout -b new-branches",
            "git create new-branches"
        ],
        "score": "git checkout -b new-branches"
    },
    {
        "title": "True or False: Git\u2019s only function is as a version control tool for code.",
        "ques_type": 11,
        "options": [
            "true",
            "false"
        ],
        "score": "False"
    },
    {
        "title": "You want to obtain an update from the remote repository origin/master. Which command(s) below achieve(s) this?",
        "ques_type": 15,
        "options": [
            "git pull origin master",
            "git clone origin master",
            "git fetch origin master &amp git merge origin master",
            "git cherry pick origin master",
            "git diff origin master"
        ],
        "score": [
            "git pull origin master",
            "git fetch origin master &amp git merge origin master"
        ]
    }
]
'@

# The old layout had the text in A2 (plain row below a placeholder A1=0,
# with A1 bold/bordered/centered). The new layout has only one row: the
# text itself sits in A1 with the workbook's default (unstyled) formatting.
$ws.Rows.Item(2).Delete()

$a1 = $ws.Range("A1")
$a1.ClearFormats()
$a1.Value = $newText

# Setting a long value can make Excel auto-expand the row height; restore
# the natural/default row height so it is not pinned to a custom value.
$ws.Rows.Item(1).AutoFit()
